$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.198.40'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '3.755.34'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '3.754.29'
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.541'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000249'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").Value = '4.386.65'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").Value = '3.751.27'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").Value = '69.195.89'
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +18.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("E24").Value = '  +5.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("E31").Value = '  +4.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.62'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").Value = '3.901.34'
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = '3.690.66'
$ws.Range("E36").Value = '  +0.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("E40").Value = '  +0.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.325'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.83%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '433.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("D50").Value = '2.791.02'
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0353'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '
